# proportion_calculations.csv: replace the hard-coded "% unique" summary
# values in the B51:D55 block with live formulas that pull from the
# per-category J-column totals computed earlier on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51 - Tornado
$ws.Range("B51").Formula = "=J3"

# Row 52 - Hail
$ws.Range("B52").Formula = "=J13"

# Row 53 - Wind
$ws.Range("B53").Formula = "=J22"
$ws.Range("C53").Formula = "=J25"
$ws.Range("D53").Formula = "=J28"

# Row 54 - Flood
$ws.Range("B54").Formula = "=J32"
$ws.Range("C54").Formula = "=J35"
$ws.Range("D54").Formula = "=J38"

# Row 55 - Winter
$ws.Range("B55").Formula = "=J42"
$ws.Range("C55").Formula = "=J45"
$ws.Range("D55").Formula = "=J48"

# Move the viewport / active selection down to where the user left off
# editing (was B51, now D56).
$ws.Range("D56").Select()
